$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (borders/styles) from the last existing data row (31) down
# across the 30 new rows (32:61) that are about to be populated.
$ws.Range("A31:D31").Copy()
$ws.Range("A32:D61").PasteSpecial(-4122)

# S.no / Power / OP values for the new rows, mirroring rows 2:31.
$data = @(
    @(150, 7654),
    @(160, 4567),
    @(170, 9876),
    @(180, 5432),
    @(190, 6543),
    @(200, 8764),
    @(210, 9876),
    @(220, 3456),
    @(230, 6545),
    @(240, 1234),
    @(250, 6432),
    @(260, 9786),
    @(270, 8977),
    @(280, 3456),
    @(290, 7866),
    @(300, 5353),
    @(310, 7564),
    @(320, 4738),
    @(330, 2468),
    @(340, 8346),
    @(341, 8347),
    @(342, 8348),
    @(343, 8349),
    @(344, 8350),
    @(345, 8351),
    @(346, 8352),
    @(347, 8353),
    @(348, 8354),
    @(349, 8355),
    @(350, 8356)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = 32 + $i
    $sno = 31 + $i
    $ws.Cells.Item($rowNum, 1).Value = $sno
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][1]
}

# CT/R column: first new row gets its own formula, the remaining 29 are
# entered as one multi-cell fill (mirrors how rows 3:31 became a shared
# formula group in the original sheet).
$ws.Range("D32").Formula = "=AVERAGE(B32:C32)"
$ws.Range("D33:D61").Formula = "=AVERAGE(B33:C33)"

$ws.Range("A31:A61").Select() | Out-Null
